# Re-saving the deck a day later causes PowerPoint to recompute the
# auto-updating "datetimeFigureOut" date fields shown on the handout
# and notes masters (Insert > Header & Footer > Date and time).
# Update both masters' date placeholders from 11/26/2023 to 11/27/2023.

$p = $ppt.ActivePresentation

$handoutDate = $p.HandoutMaster.HeadersFooters.DateAndTime
$handoutDate.Text = "11/27/2023"

$notesDate = $p.NotesMaster.HeadersFooters.DateAndTime
$notesDate.Text = "11/27/2023"
